# Merge the split "word" + "space" + "word" runs in the Title, Author and
# Date paragraphs into a single run each, matching the target edit:
#   "VIP" + " " + "report"   -> "VIP report"
#   "Tom" + " " + "Coleman"  -> "Tom Coleman"
#   "Invalid" + " " + "Date" -> "Invalid Date"
#
# Using Find/Replace on the whole (already-correct) text causes Word to
# collapse the paragraph's multiple runs into the single run that now holds
# the replacement text, which is exactly the structural change shown in the
# diff.

$d = $word.ActiveDocument

$d.Content.Find.Execute("VIP report", $true, $false, $false, $false, $false,
                         $true, 1, $false, "VIP report", 2)

$d.Content.Find.Execute("Tom Coleman", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Tom Coleman", 2)

$d.Content.Find.Execute("Invalid Date", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Invalid Date", 2)
